# Updates cryptos list cell values (Price/Volume columns) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting numeric-looking strings (e.g. "7.38") into real numbers.
# The cells original style is preserved (no lasting format change).
function Set-TextValue($range, $text) {
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $originalStyle
}

$ws.Range("D2").Value = '27.772.49'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '1.595.34'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '209.28'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("E9").Value = '  -1.32%  '
Set-TextValue $ws.Range("D10") '0.0593'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '1.823.04'
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '1.603.18'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("E14").Value = '  -2.20%  '
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").Value = '27.776.21'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("E17").Value = '  -1.45%  '
Set-TextValue $ws.Range("D18") '219.63'
$ws.Range("E18").Value = '  -2.72%  '
$ws.Range("D19").Value = '0.0₃0697'
$ws.Range("E19").Value = '  -2.05%  '
Set-TextValue $ws.Range("D20") '7.38'
$ws.Range("E20").Value = '  -2.34%  '
$ws.Range("E22").Value = '  -3.22%  '
Set-TextValue $ws.Range("D23") '9.76'
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("E24").Value = '  -3.79%  '
Set-TextValue $ws.Range("D25") '154.03'
$ws.Range("E25").Value = '  -0.04%  '
Set-TextValue $ws.Range("D26") '7.23'
$ws.Range("E26").Value = '  +4.97%  '
$ws.Range("E27").Value = '  +0.03%  '
Set-TextValue $ws.Range("D29") '0.106'
$ws.Range("E29").Value = '  -3.59%  '
Set-TextValue $ws.Range("D30") '1.17'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D33").Value = '1.377.26'
$ws.Range("E33").Value = '  -2.65%  '
Set-TextValue $ws.Range("D34") '2.97'
$ws.Range("E34").Value = '  -2.70%  '
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  +0.05%  '
Set-TextValue $ws.Range("D39") '0.537'
$ws.Range("E39").Value = '  -2.58%  '
Set-TextValue $ws.Range("D40") '0.830'
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("E42").Value = '  -2.71%  '
Set-TextValue $ws.Range("D43") '64.58'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("E44").Value = '  +2.43%  '
Set-TextValue $ws.Range("D45") '5.25'
$ws.Range("E45").Value = '  -1.68%  '
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").Value = '1.733.99'
$ws.Range("E47").Value = '  -1.60%  '
Set-TextValue $ws.Range("D48") '86.63'
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("E49").Value = '  -0.29%  '
Set-TextValue $ws.Range("D50") '0.0967'
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("E51").Value = '  -0.81%  '
